# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the header row's "_old"/"_new" suffixes to the respective
# format-version suffixes ("_FV2410" / "_FV2504"), freezes the header row,
# and (re-)wraps the used range in an Excel Table ("Table1") so the
# AutoFilter dropdowns on the header row keep working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row -------------------------------------------------
$oldSuffixCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newSuffixCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldSuffixCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
    $ws.Range($newSuffixCols[$i] + "1").Value = $baseNames[$i] + "_FV2504"
}

# K1 ("diff") is unchanged.

# --- 2. Freeze the header row ---------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap the used range in a table so the header row keeps its filter --
$usedRange = $ws.Range("A1:U59")
$table = $ws.ListObjects.Add(1, $usedRange, [Type]::Missing, 1)
$table.Name = "Table1"
